$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the duplicated, bold "Play Coywolf Cash Free..." paragraph near
#    the bottom of the document (the title already appears as the Heading1
#    at the top; this one at the bottom is redundant and gets deleted).
#    Doing this first (before the text further down is duplicated at the
#    top of the doc) keeps the later text search unambiguous.
# ---------------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^Play Coywolf Cash Free: Immersive American Wilderness Slot Game\r?$" -and $p.Range.Font.Bold) {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Replace the text of the remaining italic "meta description" paragraph
#    at the bottom with the new image-prompt text, keeping its formatting.
# ---------------------------------------------------------------------------
$newBottomText = "Please create a cartoon image for Coywolf Cash featuring a happy Maya warrior with glasses. The image should be fun and engaging, with bright colors and cartoon-style graphics. The Maya warrior should be smiling and holding a bag of money adorned with a dollar sign, with the Coywolf Cash slot machine in the background. The background of the image should feature the American wilderness, with rock formations, cacti and the endless road stretching out into the distance. The image should be eye-catching and encourage potential players to give Coywolf Cash a try."

$null = $d.Content.Find.Execute("Read our game review for Coywolf Cash and play for free. Enjoy an immersive, American wilderness themed slot game with a high payout potential.", $true, $false, $false, $false, $false, $true, 1, $false, $newBottomText, 2)

# ---------------------------------------------------------------------------
# 3. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)            # collapse to the end of the title paragraph
$titleRange.InsertParagraphAfter() # creates a new, empty paragraph after it

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = $d.Styles.Item("Normal")   # make sure it is a normal body paragraph

$metaRange = $metaPara.Range
$metaRange.Collapse(1)             # move to the very start of the new paragraph

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our game review for Coywolf Cash and play for free. Enjoy an immersive, American wilderness themed slot game with a high payout potential.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $metaRange.InsertXML($metaXml)
